# Apply cryptocurrency price/volume updates to the sheet.
# Source data is plain text in column D (price) and E (1h volume %),
# with B/C (name/link) swapped for a few re-ranked coins.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '43.708.70'
$ws.Range("E2").Value = '  +1.85%  '

# Row 3
$ws.Range("D3").Value = '2.209.47'
$ws.Range("E3").Value = '  -0.45%  '

# Row 4
$ws.Range("E4").Value = '  +0.15%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '265.85'
$ws.Range("E5").Value = '  +3.33%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '86.08'
$ws.Range("E6").Value = '  +11.58%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.620'
$ws.Range("E7").Value = '  +0.16%  '

# Row 8
$ws.Range("E8").Value = '  +0.07%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.603'
$ws.Range("E9").Value = '  +0.98%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '46.03'
$ws.Range("E10").Value = '  +8.76%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0918'
$ws.Range("E11").Value = '  +1.24%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.52'
$ws.Range("E12").Value = '  +7.30%  '

# Row 13
$ws.Range("E13").Value = '  +1.98%  '

# Row 14
$ws.Range("D14").Value = '2.541.83'
$ws.Range("E14").Value = '  -0.22%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.59'
$ws.Range("E15").Value = '  +0.43%  '

# Row 16
$ws.Range("D16").Value = '2.226.09'
$ws.Range("E16").Value = '  +0.55%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.783'
$ws.Range("E17").Value = '  -0.25%  '

# Row 18
$ws.Range("D18").Value = '43.664.96'
$ws.Range("E18").Value = '  +1.79%  '

# Row 19
$ws.Range("E19").Value = '  +0.36%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.97'
$ws.Range("E20").Value = '  -0.43%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '69.92'
$ws.Range("E21").Value = '  -1.81%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.36'
$ws.Range("E22").Value = '  +6.50%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '231.53'
$ws.Range("E23").Value = '  +0.27%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.81'
$ws.Range("E24").Value = '  -5.59%  '

# Row 25
$ws.Range("B25").Value = 'PancakeSwap'
$ws.Range("C25").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.65'
$ws.Range("E25").Value = '  +20.39%  '

# Row 26
$ws.Range("B26").Value = 'Dai'
$ws.Range("C26").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("E26").Value = '  -0.04%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.78'
$ws.Range("E27").Value = '  -0.26%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.54'
$ws.Range("E28").Value = '  +6.36%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '39.14'
$ws.Range("E29").Value = '  -9.27%  '

# Row 30
$ws.Range("E30").Value = '  +1.75%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '174.79'
$ws.Range("E31").Value = '  +1.02%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0890'
$ws.Range("E32").Value = '  +1.40%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '20.49'
$ws.Range("E33").Value = '  +0.17%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.38'
$ws.Range("E34").Value = '  +2.62%  '

# Row 35
$ws.Range("E35").Value = '  +1.40%  '

# Row 36
$ws.Range("E36").Value = '  +1.29%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0356'
$ws.Range("E37").Value = '  -1.84%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.38'
$ws.Range("E38").Value = '  +0.89%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.27'
$ws.Range("E39").Value = '  +16.05%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '12.24'
$ws.Range("E40").Value = '  -5.35%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '64.59'
$ws.Range("E41").Value = '  +7.14%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.08'
$ws.Range("E42").Value = '  -1.45%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.46'
$ws.Range("E43").Value = '  +2.88%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.202'
$ws.Range("E44").Value = '  -0.25%  '

# Row 45
$ws.Range("B45").Value = 'Cronos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0985'
$ws.Range("E45").Value = '  +0.58%  '

# Row 46
$ws.Range("B46").Value = 'FraxShare'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.32'
$ws.Range("E46").Value = '  -0.44%  '

# Row 47
$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '100.06'
$ws.Range("E47").Value = '  -2.92%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.18'
$ws.Range("E48").Value = '  +3.99%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.12'
$ws.Range("E49").Value = '  -0.17%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.440'
$ws.Range("E50").Value = '  -5.61%  '

# Row 51
$ws.Range("E51").Value = '  +4.74%  '
